$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new date cells (rows 19-29, column D) use the same date number format
# as the rest of the "Fecha" column before assigning their values.
$ws.Range("D19:D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---- Update existing rows 5-18 (changed cells only) ----
# Row 5
$ws.Range("D5").Value = 44425
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 35
$ws.Range("N5").Value = 4500
$ws.Range("O5").Value = 4500
$ws.Range("P5").Value = 4500
$ws.Range("S5").Value = 4500

# Row 6
$ws.Range("D6").Value = 44425
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 20

# Row 7
$ws.Range("D7").Value = 44425
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("S7").Value = 3000

# Row 8
$ws.Range("D8").Value = 44411
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("Q8").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 1

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("S9").Value = 3000

# Row 10
$ws.Range("D10").Value = 44407
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 40

# Row 11
$ws.Range("D11").Value = 44414
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 55

# Row 12
$ws.Range("D12").Value = 44162
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = 2200
$ws.Range("O12").Value = 2300
$ws.Range("P12").Value = 2247
$ws.Range("S12").Value = 2247

# Row 13
$ws.Range("D13").Value = 44427
$ws.Range("L13").Value = 'Especial'
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 24000
$ws.Range("Q13").Value = '$/bandeja 7 kilos'
$ws.Range("S13").Value = 3429
$ws.Range("T13").Value = 7

# Row 14
$ws.Range("D14").Value = 44413
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 35
$ws.Range("N14").Value = 3500
$ws.Range("O14").Value = 3500
$ws.Range("P14").Value = 3500
$ws.Range("S14").Value = 3500

# Row 15
$ws.Range("D15").Value = 44421
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 3200
$ws.Range("O15").Value = 3200
$ws.Range("P15").Value = 3200
$ws.Range("S15").Value = 3200

# Row 16
$ws.Range("D16").Value = 44377
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 3500
$ws.Range("O16").Value = 3500
$ws.Range("P16").Value = 3500
$ws.Range("S16").Value = 3500

# Row 17
$ws.Range("D17").Value = 44426
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 4500
$ws.Range("O17").Value = 4500
$ws.Range("P17").Value = 4500
$ws.Range("S17").Value = 4500

# Row 18
$ws.Range("D18").Value = 44426
$ws.Range("M18").Value = 45
$ws.Range("N18").Value = 3500
$ws.Range("O18").Value = 3500
$ws.Range("P18").Value = 3500
$ws.Range("S18").Value = 3500

# ---- Add new rows 19-29 (full rows) ----
# Row 19
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = 'Vega Modelo de Temuco'
$ws.Range("C19").Value = 'La Araucanía'
$ws.Range("D19").Value = 44354
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 'Fruta'
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = 'Otros'
$ws.Range("I19").Value = 100107002
$ws.Range("J19").Value = 'Chirimoya'
$ws.Range("K19").Value = 'Cultivar IV Región'
$ws.Range("L19").Value = 'Tercera'
$ws.Range("M19").Value = 95
$ws.Range("N19").Value = 3500
$ws.Range("O19").Value = 3500
$ws.Range("P19").Value = 3500
$ws.Range("Q19").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R19").Value = 'Provincia del Elquí'
$ws.Range("S19").Value = 3500
$ws.Range("T19").Value = 1

# Row 20
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = 'Vega Modelo de Temuco'
$ws.Range("C20").Value = 'La Araucanía'
$ws.Range("D20").Value = 44412
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 'Fruta'
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = 'Otros'
$ws.Range("I20").Value = 100107002
$ws.Range("J20").Value = 'Chirimoya'
$ws.Range("K20").Value = 'Cultivar IV Región'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 65
$ws.Range("N20").Value = 3200
$ws.Range("O20").Value = 3200
$ws.Range("P20").Value = 3200
$ws.Range("Q20").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R20").Value = 'Provincia del Elquí'
$ws.Range("S20").Value = 3200
$ws.Range("T20").Value = 1

# Row 21
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = 'Vega Modelo de Temuco'
$ws.Range("C21").Value = 'La Araucanía'
$ws.Range("D21").Value = 44405
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 'Fruta'
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = 'Otros'
$ws.Range("I21").Value = 100107002
$ws.Range("J21").Value = 'Chirimoya'
$ws.Range("K21").Value = 'Cultivar IV Región'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 3200
$ws.Range("O21").Value = 3200
$ws.Range("P21").Value = 3200
$ws.Range("Q21").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R21").Value = 'Provincia del Elquí'
$ws.Range("S21").Value = 3200
$ws.Range("T21").Value = 1

# Row 22
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = 'Vega Modelo de Temuco'
$ws.Range("C22").Value = 'La Araucanía'
$ws.Range("D22").Value = 44417
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 'Fruta'
$ws.Range("G22").Value = 100107
$ws.Range("H22").Value = 'Otros'
$ws.Range("I22").Value = 100107002
$ws.Range("J22").Value = 'Chirimoya'
$ws.Range("K22").Value = 'Cultivar IV Región'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 3200
$ws.Range("O22").Value = 3200
$ws.Range("P22").Value = 3200
$ws.Range("Q22").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R22").Value = 'Provincia del Elquí'
$ws.Range("S22").Value = 3200
$ws.Range("T22").Value = 1

# Row 23
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = 'Vega Modelo de Temuco'
$ws.Range("C23").Value = 'La Araucanía'
$ws.Range("D23").Value = 44419
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 'Fruta'
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = 'Otros'
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = 'Chirimoya'
$ws.Range("K23").Value = 'Cultivar IV Región'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 70
$ws.Range("N23").Value = 3200
$ws.Range("O23").Value = 3200
$ws.Range("P23").Value = 3200
$ws.Range("Q23").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R23").Value = 'Provincia del Elquí'
$ws.Range("S23").Value = 3200
$ws.Range("T23").Value = 1

# Row 24
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 'Vega Modelo de Temuco'
$ws.Range("C24").Value = 'La Araucanía'
$ws.Range("D24").Value = 44420
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 'Fruta'
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = 'Otros'
$ws.Range("I24").Value = 100107002
$ws.Range("J24").Value = 'Chirimoya'
$ws.Range("K24").Value = 'Cultivar IV Región'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 35
$ws.Range("N24").Value = 3500
$ws.Range("O24").Value = 3500
$ws.Range("P24").Value = 3500
$ws.Range("Q24").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R24").Value = 'Provincia del Elquí'
$ws.Range("S24").Value = 3500
$ws.Range("T24").Value = 1

# Row 25
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = 'Vega Modelo de Temuco'
$ws.Range("C25").Value = 'La Araucanía'
$ws.Range("D25").Value = 44420
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 'Fruta'
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = 'Otros'
$ws.Range("I25").Value = 100107002
$ws.Range("J25").Value = 'Chirimoya'
$ws.Range("K25").Value = 'Cultivar IV Región'
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 3200
$ws.Range("O25").Value = 3200
$ws.Range("P25").Value = 3200
$ws.Range("Q25").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R25").Value = 'Provincia del Elquí'
$ws.Range("S25").Value = 3200
$ws.Range("T25").Value = 1

# Row 26
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = 'Vega Modelo de Temuco'
$ws.Range("C26").Value = 'La Araucanía'
$ws.Range("D26").Value = 44161
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 'Fruta'
$ws.Range("G26").Value = 100107
$ws.Range("H26").Value = 'Otros'
$ws.Range("I26").Value = 100107002
$ws.Range("J26").Value = 'Chirimoya'
$ws.Range("K26").Value = 'Cultivar IV Región'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 65
$ws.Range("N26").Value = 2300
$ws.Range("O26").Value = 2300
$ws.Range("P26").Value = 2300
$ws.Range("Q26").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R26").Value = 'Provincia del Elquí'
$ws.Range("S26").Value = 2300
$ws.Range("T26").Value = 1

# Row 27
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = 'Vega Modelo de Temuco'
$ws.Range("C27").Value = 'La Araucanía'
$ws.Range("D27").Value = 44161
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = 'Otros'
$ws.Range("I27").Value = 100107002
$ws.Range("J27").Value = 'Chirimoya'
$ws.Range("K27").Value = 'Cultivar IV Región'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 55
$ws.Range("N27").Value = 2000
$ws.Range("O27").Value = 2000
$ws.Range("P27").Value = 2000
$ws.Range("Q27").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R27").Value = 'Provincia del Elquí'
$ws.Range("S27").Value = 2000
$ws.Range("T27").Value = 1

# Row 28
$ws.Range("A28").Value = 10
$ws.Range("B28").Value = 'Vega Modelo de Temuco'
$ws.Range("C28").Value = 'La Araucanía'
$ws.Range("D28").Value = 44159
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 'Fruta'
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = 'Otros'
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = 'Chirimoya'
$ws.Range("K28").Value = 'Cultivar IV Región'
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 2300
$ws.Range("O28").Value = 2500
$ws.Range("P28").Value = 2408
$ws.Range("Q28").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R28").Value = 'Provincia del Elquí'
$ws.Range("S28").Value = 2408
$ws.Range("T28").Value = 1

# Row 29
$ws.Range("A29").Value = 10
$ws.Range("B29").Value = 'Vega Modelo de Temuco'
$ws.Range("C29").Value = 'La Araucanía'
$ws.Range("D29").Value = 44160
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 'Fruta'
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = 'Otros'
$ws.Range("I29").Value = 100107002
$ws.Range("J29").Value = 'Chirimoya'
$ws.Range("K29").Value = 'Cultivar IV Región'
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 120
$ws.Range("N29").Value = 2200
$ws.Range("O29").Value = 2300
$ws.Range("P29").Value = 2246
$ws.Range("Q29").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R29").Value = 'Provincia del Elquí'
$ws.Range("S29").Value = 2246
$ws.Range("T29").Value = 1

